$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: fill in previously empty cells C10:F10 with 5
$ws.Range("C10:F10").Value = 5

# Row 23: fill in previously empty cells F23:G23 with 5 (H23 stays empty)
$ws.Range("F23:G23").Value = 5

# Move/record the active selection on the active pane to H23
$ws.Range("H23").Select()
